$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "efef efe" "Tovstukha Eduard"
Replace-Text "frfd rfrd" "Molodiznaa 12/12q"
Replace-Text "434 efef" "32233 Chmelnitskiy"
Replace-Text "dfrf" "Ukraine"
Replace-Text "de@efe" "ed@ed.com"
Replace-Text "4334" "380985351072"
Replace-Text "drfdrf" "swswsw"
Replace-Text "Amount USD: 254" "Amount USD: 381"
Replace-Text "___________________26-3-2020" "___________________18-4-2020"
